$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(1,350).Value = 0.1
$ws.Cells.Item(1,351).Value = 1.5
$ws.Cells.Item(1,352).Value = 100.41
$ws.Cells.Item(1,353).Value = 2.0
$ws.Cells.Item(1,354).Value = 3.25
